$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 103.875
$ws.Range("I55").Value = 111.57143
$ws.Range("J55").Value = 50
$ws.Range("K55").Value = 111.57143
$ws.Range("L55").Value = 50
$ws.Range("M55").Value = 102.42857
$ws.Range("N55").Value = -478

$ws.Range("H98").Value = 1100.6875
$ws.Range("I98").Value = 907.2
$ws.Range("J98").Value = 4003
$ws.Range("K98").Value = 907.2
$ws.Range("L98").Value = 4003
$ws.Range("M98").Value = 590.8
$ws.Range("N98").Value = -6999

$ws.Range("H112").Value = 1465.4546
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 1465.4546
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 4396.3638
$ws.Range("N112").Value = -6612.3638
$ws.Range("M112").ClearContents()

$ws.Range("H122").Value = 1100.6875
$ws.Range("I122").Value = 907.2
$ws.Range("J122").Value = 4003
$ws.Range("K122").Value = 2721.6
$ws.Range("L122").Value = 12009
$ws.Range("M122").Value = -271.6000000000004
$ws.Range("N122").Value = -16909

$ws.Range("H123").Value = 30333.334
$ws.Range("J123").Value = 30333.334
$ws.Range("L123").Value = 30333.334
$ws.Range("N123").Value = -40133.334

$ws.Range("H132").Value = 3257.3052
$ws.Range("I132").Value = 1582.5962
$ws.Range("K132").Value = 4747.7886
$ws.Range("M132").Value = -2217.7886

$ws.Range("H137").Value = 1324177.6
$ws.Range("I137").Value = 1532.0646
$ws.Range("J137").Value = 5051633.5
$ws.Range("K137").Value = 4596.1938
$ws.Range("L137").Value = 15154900.5
$ws.Range("M137").Value = -2046.1938
$ws.Range("N137").Value = -15160000.5

$ws.Range("H138").Value = 3127396.2
$ws.Range("I138").Value = 1055.1945
$ws.Range("J138").Value = 7146978
$ws.Range("K138").Value = 3165.5835
$ws.Range("L138").Value = 21440934
$ws.Range("M138").Value = 1974.4165
$ws.Range("N138").Value = -21451214

$ws.Range("H141").Value = 828.902
$ws.Range("I141").Value = 774.4681
$ws.Range("J141").Value = 1468.5
$ws.Range("K141").Value = 2323.4043
$ws.Range("L141").Value = 4405.5
$ws.Range("M141").Value = 2856.5957
$ws.Range("N141").Value = -14765.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1068.9639
$ws.Range("I61").Value = 926.7222
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 926.7222
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -714.7222
$ws.Range("N61").Value = -2424

$ws.Range("H74").Value = 22802.176
$ws.Range("I74").Value = 26489.871
$ws.Range("J74").Value = 10817.167
$ws.Range("K74").Value = 26489.871
$ws.Range("L74").Value = 10817.167
$ws.Range("M74").Value = -25615.871
$ws.Range("N74").Value = -12565.167

$ws.Range("H77").Value = 22802.176
$ws.Range("I77").Value = 26489.871
$ws.Range("J77").Value = 10817.167
$ws.Range("K77").Value = 132449.355
$ws.Range("L77").Value = 54085.835
$ws.Range("M77").Value = -128081.355
$ws.Range("N77").Value = -62821.835

$ws.Range("H136").Value = 1068.9639
$ws.Range("I136").Value = 926.7222
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 2780.1666
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -230.1666
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1879.8
$ws.Range("J86").Value = 2750
$ws.Range("L86").Value = 2750
$ws.Range("N86").Value = -4996

$ws.Range("H89").Value = 1879.8
$ws.Range("J89").Value = 2750
$ws.Range("L89").Value = 13750
$ws.Range("N89").Value = -24982

$ws.Range("H134").Value = 957896.75
$ws.Range("I134").Value = 1672361
$ws.Range("J134").Value = 5277.6665
$ws.Range("K134").Value = 5017083
$ws.Range("L134").Value = 15832.9995
$ws.Range("M134").Value = -5014548
$ws.Range("N134").Value = -20902.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9092012
$ws.Range("I31").Value = 1026.6522
$ws.Range("J31").Value = 55557050
$ws.Range("K31").Value = 1026.6522
$ws.Range("L31").Value = 55557050
$ws.Range("M31").Value = -731.6522
$ws.Range("N31").Value = -55557640

$ws.Range("H34").Value = 9092012
$ws.Range("I34").Value = 1026.6522
$ws.Range("J34").Value = 55557050
$ws.Range("K34").Value = 1026.6522
$ws.Range("L34").Value = 55557050
$ws.Range("M34").Value = -824.6522
$ws.Range("N34").Value = -55557454

$ws.Range("H58").Value = 3970.5588
$ws.Range("I58").Value = 4148.355
$ws.Range("J58").Value = 2133.3333
$ws.Range("K58").Value = 4148.355
$ws.Range("L58").Value = 2133.3333
$ws.Range("M58").Value = -3945.355
$ws.Range("N58").Value = -2539.3333

$ws.Range("H86").Value = 7446.7
$ws.Range("I86").Value = 6921.0386
$ws.Range("J86").Value = 8422.929
$ws.Range("K86").Value = 6921.0386
$ws.Range("L86").Value = 8422.929
$ws.Range("M86").Value = -5798.0386
$ws.Range("N86").Value = -10668.929

$ws.Range("H89").Value = 7446.7
$ws.Range("I89").Value = 6921.0386
$ws.Range("J89").Value = 8422.929
$ws.Range("K89").Value = 34605.193
$ws.Range("L89").Value = 42114.645
$ws.Range("M89").Value = -28989.193
$ws.Range("N89").Value = -53346.645

$ws.Range("H132").Value = 554205.9399999999
$ws.Range("I132").Value = 1491.2858
$ws.Range("J132").Value = 3368026
$ws.Range("K132").Value = 4473.857400000001
$ws.Range("L132").Value = 10104078
$ws.Range("M132").Value = -1943.857400000001
$ws.Range("N132").Value = -10109138

$ws.Range("H134").Value = 2225.705
$ws.Range("I134").Value = 2148.8833
$ws.Range("J134").Value = 2481.7778
$ws.Range("K134").Value = 6446.6499
$ws.Range("L134").Value = 7445.3334
$ws.Range("M134").Value = -3911.6499
$ws.Range("N134").Value = -12515.3334

$ws.Range("H136").Value = 3970.5588
$ws.Range("I136").Value = 4148.355
$ws.Range("J136").Value = 2133.3333
$ws.Range("K136").Value = 12445.065
$ws.Range("L136").Value = 6399.999899999999
$ws.Range("M136").Value = -9895.064999999999
$ws.Range("N136").Value = -11499.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 366632.97
$ws.Range("I2").Value = 1137.5555
$ws.Range("J2").Value = 585930.2
$ws.Range("K2").Value = 6825.333
$ws.Range("L2").Value = 3515581.2
$ws.Range("M2").Value = -6712.333
$ws.Range("N2").Value = -3515807.2

$ws.Range("H5").Value = 1353.25
$ws.Range("J5").Value = 1515
$ws.Range("L5").Value = 4545
$ws.Range("N5").Value = -4769

$ws.Range("H98").Value = 558
$ws.Range("I98").Value = 395.66666
$ws.Range("J98").Value = 923.25
$ws.Range("K98").Value = 1186.99998
$ws.Range("L98").Value = 2769.75
$ws.Range("M98").Value = 311.0000199999999
$ws.Range("N98").Value = -5765.75

$ws.Range("H123").Value = 5588.8887
$ws.Range("J123").Value = 6614.2856
$ws.Range("L123").Value = 19842.8568
$ws.Range("N123").Value = -24742.8568

$ws.Range("H129").Value = 2244.1428
$ws.Range("I129").Value = 3137
$ws.Range("J129").Value = 1748.1111
$ws.Range("K129").Value = 9411
$ws.Range("L129").Value = 5244.3333
$ws.Range("M129").Value = -4411
$ws.Range("N129").Value = -15244.3333

$ws.Range("H131").Value = 937.22
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 937.22
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 2811.66
$ws.Range("N131").Value = -12891.66
$ws.Range("M131").ClearContents()

$ws.Range("H133").Value = 6865.35
$ws.Range("J133").Value = 7031.25
$ws.Range("L133").Value = 21093.75
$ws.Range("N133").Value = -31213.75

$ws.Range("H134").Value = 9792.666999999999
$ws.Range("I134").Value = 10270
$ws.Range("J134").Value = 9375
$ws.Range("K134").Value = 30810
$ws.Range("L134").Value = 28125
$ws.Range("M134").Value = -25740
$ws.Range("N134").Value = -38265

$ws.Range("H135").Value = 1353.25
$ws.Range("J135").Value = 1515
$ws.Range("L135").Value = 13635
$ws.Range("N135").Value = -18705

$ws.Range("H136").Value = 4349
$ws.Range("I136").Value = 1497.5
$ws.Range("J136").Value = 6250
$ws.Range("K136").Value = 4492.5
$ws.Range("L136").Value = 18750
$ws.Range("M136").Value = 607.5
$ws.Range("N136").Value = -28950

$ws.Range("H137").Value = 16512016
$ws.Range("I137").Value = 3618
$ws.Range("J137").Value = 18804850
$ws.Range("K137").Value = 10854
$ws.Range("L137").Value = 56414550
$ws.Range("M137").Value = -5754
$ws.Range("N137").Value = -56424750

$ws.Range("H138").Value = 1402.6471
$ws.Range("I138").Value = 1119.9286
$ws.Range("K138").Value = 3359.7858
$ws.Range("M138").Value = 1780.2142

$ws.Range("H139").Value = 1874.8948
$ws.Range("I139").Value = 1235.4546
$ws.Range("K139").Value = 3706.3638
$ws.Range("M139").Value = 1433.6362

$ws.Range("H140").Value = 2824.3635
$ws.Range("I140").Value = 906.8
$ws.Range("K140").Value = 2720.4
$ws.Range("M140").Value = 2459.6

$ws.Range("H141").Value = 2576
$ws.Range("I141").Value = 2195.5557
$ws.Range("K141").Value = 6586.6671
$ws.Range("M141").Value = -1406.6671

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3641.6724
$ws.Range("I132").Value = 3255.087
$ws.Range("J132").Value = 5123.5835
$ws.Range("K132").Value = 9765.261
$ws.Range("L132").Value = 15370.7505
$ws.Range("M132").Value = -7235.261
$ws.Range("N132").Value = -20430.7505

$ws.Range("H136").Value = 2134.9395
$ws.Range("I136").Value = 1330.619
$ws.Range("J136").Value = 3542.5
$ws.Range("K136").Value = 3991.857
$ws.Range("L136").Value = 10627.5
$ws.Range("M136").Value = -1441.857
$ws.Range("N136").Value = -15727.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2698.5112
$ws.Range("I132").Value = 2765.6765
$ws.Range("J132").Value = 2490.9092
$ws.Range("K132").Value = 8297.029500000001
$ws.Range("L132").Value = 7472.7276
$ws.Range("M132").Value = -5767.029500000001
$ws.Range("N132").Value = -12532.7276

$ws.Range("H136").Value = 1458.45
$ws.Range("I136").Value = 1316.7742
$ws.Range("J136").Value = 1946.4445
$ws.Range("K136").Value = 3950.3226
$ws.Range("L136").Value = 5839.333500000001
$ws.Range("M136").Value = -1400.3226
$ws.Range("N136").Value = -10939.3335
